# Apply cryptocurrency price/volume update scraped on Mon Jan  1 03:39:30 UTC 2024
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D and E to be treated as text so that numeric-looking values
# (e.g. "8.20", "0.0900", "42.426.13") are preserved exactly as strings instead
# of being converted into numbers by Excel.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "42.426.13"
$ws.Range("E2").Value = "  +0.77%  "
# Row 3
$ws.Range("D3").Value = "2.282.50"
$ws.Range("E3").Value = "  +0.09%  "
# Row 4
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.35%  "
# Row 5
$ws.Range("D5").Value = "310.80"
$ws.Range("E5").Value = "  -2.62%  "
# Row 6
$ws.Range("D6").Value = "103.70"
$ws.Range("E6").Value = "  +3.20%  "
# Row 7
$ws.Range("E7").Value = "  +0.33%  "
# Row 8
$ws.Range("E8").Value = "  -0.03%  "
# Row 9
$ws.Range("D9").Value = "0.598"
$ws.Range("E9").Value = "  -0.35%  "
# Row 10
$ws.Range("D10").Value = "38.84"
$ws.Range("E10").Value = "  +0.08%  "
# Row 11
$ws.Range("D11").Value = "0.0900"
$ws.Range("E11").Value = "  +0.19%  "
# Row 12
$ws.Range("D12").Value = "8.20"
$ws.Range("E12").Value = "  -0.11%  "
# Row 13
$ws.Range("E13").Value = "  +1.75%  "
# Row 14
$ws.Range("E14").Value = "  +2.59%  "
# Row 15
$ws.Range("D15").Value = "15.03"
$ws.Range("E15").Value = "  +0.12%  "
# Row 16
$ws.Range("D16").Value = "2.628.40"
$ws.Range("E16").Value = "  +0.08%  "
# Row 17
$ws.Range("D17").Value = "2.288.19"
$ws.Range("E17").Value = "  +0.79%  "
# Row 18
$ws.Range("D18").Value = "42.363.45"
$ws.Range("E18").Value = "  +0.40%  "
# Row 19
$ws.Range("E19").Value = "  -0.47%  "
# Row 20
$ws.Range("E20").Value = "  -0.14%  "
# Row 21
$ws.Range("D21").Value = "13.46"
$ws.Range("E21").Value = "  +6.43%  "
# Row 22
$ws.Range("D22").Value = "72.81"
$ws.Range("E22").Value = "  +0.41%  "
# Row 23
$ws.Range("E23").Value = "  -3.26%  "
# Row 24
$ws.Range("D24").Value = "262.90"
$ws.Range("E24").Value = "  -1.64%  "
# Row 25
$ws.Range("E25").Value = "  -1.26%  "
# Row 26
$ws.Range("E26").Value = "  +0.34%  "
# Row 27
$ws.Range("E27").Value = "  -0.83%  "
# Row 28
$ws.Range("E28").Value = "  +16.17%  "
# Row 29
$ws.Range("D29").Value = "2.26"
$ws.Range("E29").Value = "  -2.52%  "
# Row 30
$ws.Range("D30").Value = "22.21"
$ws.Range("E30").Value = "  -0.51%  "
# Row 31
$ws.Range("D31").Value = "35.69"
$ws.Range("E31").Value = "  -3.99%  "
# Row 32
$ws.Range("D32").Value = "164.47"
$ws.Range("E32").Value = "  +0.45%  "
# Row 33
$ws.Range("E33").Value = "  -0.93%  "
# Row 34
$ws.Range("E34").Value = "  -1.76%  "
# Row 35
$ws.Range("D35").Value = "2.56"
$ws.Range("E35").Value = "  +1.60%  "
# Row 36
$ws.Range("E36").Value = "  -2.53%  "
# Row 38
$ws.Range("E38").Value = "  -1.33%  "
# Row 39
$ws.Range("D39").Value = "3.72"
$ws.Range("E39").Value = "  +2.07%  "
# Row 40
$ws.Range("D40").Value = "2.71"
$ws.Range("E40").Value = "  -1.85%  "
# Row 41
$ws.Range("E41").Value = "  +3.08%  "
# Row 42
$ws.Range("D42").Value = "98.33"
$ws.Range("E42").Value = "  +7.25%  "
# Row 43
$ws.Range("D43").Value = "68.80"
$ws.Range("E43").Value = "  +1.30%  "
# Row 44
$ws.Range("B44").Value = "Algorand"
$ws.Range("C44").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D44").Value = "0.225"
$ws.Range("E44").Value = "  +1.20%  "
# Row 45
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").Value = "1.00"
$ws.Range("E45").Value = "  -0.40%  "
# Row 46
$ws.Range("D46").Value = "1.720.52"
$ws.Range("E46").Value = "  +7.23%  "
# Row 47
$ws.Range("D47").Value = "11.89"
$ws.Range("E47").Value = "  +0.32%  "
# Row 48
$ws.Range("D48").Value = "109.96"
$ws.Range("E48").Value = "  -4.35%  "
# Row 49
$ws.Range("D49").Value = "76.98"
$ws.Range("E49").Value = "  -2.28%  "
# Row 50
$ws.Range("E50").Value = "  -0.63%  "
# Row 51
$ws.Range("D51").Value = "8.60"
$ws.Range("E51").Value = "  -3.43%  "

# Restore the default (Normal) style on the data range so that no stray
# number-format / style is left applied to the cells.
$dataRange.Style = "Normal"

